$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing elective codes for General Arts / General Science courses
$ws.Range("D12").Value = "GEO ECO EMATH GOV"
$ws.Range("D13").Value = "PHY CHE EMATH BIO EICT"

# Add a new course row for Business
$ws.Range("B14").Value = "BUS"
$ws.Range("C14").Value = "Business"
$ws.Range("D14").Value = "FA CA BM ECO EMATH"

# Move the active selection to the newly added cell
$ws.Range("D14").Select()
